$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The URL row (row 40) loses its hyperlink and becomes a plain italic
# source line like the rows around it (matching row 38/39/41's style)
# before any rows move, so the new separator rows below all inherit the
# same plain italic look instead of the hyperlink look.
$ws.Range("A40").Hyperlinks.Delete()
$ws.Range("A40").Font.Italic = $true
$ws.Range("A40").Font.Underline = $false
$ws.Range("A40").Font.Bold = $false

# Insert a blank separator row after "Source:" (row 38), before the
# "Direccion General..." line (old row 39).
$ws.Rows("39").Insert()
$ws.Range("A39").Value = ""

# Insert a blank separator row after "Direccion General..." (now row 40),
# before the URL line (now row 41).
$ws.Rows("41").Insert()
$ws.Range("A41").Value = ""

# Insert a blank separator row after the URL line (now row 42), before
# "Section 2.5" (now row 43).
$ws.Rows("43").Insert()
$ws.Range("A43").Value = ""

# The final citation line (now A48) is shortened to just "CODEMYPE",
# matching the row above it (A47).
$ws.Range("A48").Value = "CODEMYPE"
